$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultiCardLesson")

# Row 59 - section header "LessonEdit"
$ws.Range("A59").Value = "LessonEdit"

# Row 60 - data row for LessonEdit (A60 value entered last, see below)
$ws.Range("B60").Value = "Sample Content for Text Card"
$ws.Range("C60").Value = "How many colors in rainbow?"
$ws.Range("D60").Value = "Sample Content for Text Card after edit"
$ws.Range("E60").Value = "Question after edit with options"
$ws.Range("F60").Value = "C"
$ws.Range("G60").Value = "D"

# Row 61 - section header "LessonDuplicate"
$ws.Range("A61").Value = "LessonDuplicate"

# Row 62 - data row for LessonDuplicate
$ws.Range("A62").Value = "Duplicate lesson"
$ws.Range("B62").Value = "You will be creating a copy of "
$ws.Range("C62").Value = "Name the copy"
$ws.Range("D62").Value = "Copy of"
$ws.Range("E62").Value = "CANCEL"
$ws.Range("F62").Value = "SAVE"

# A60 value added last (matches original shared-string ordering)
$ws.Range("A60").Value = "Lesson_for_Checking_edit7"

# Apply styles matching the other section header rows (style index 9: yellow fill)
$ws.Range("A59").Interior.Color = $ws.Range("A56").Interior.Color
$ws.Range("A61").Interior.Color = $ws.Range("A56").Interior.Color

# Apply number format + left alignment style on D60 and E60 (matches new cellXfs index 15)
$ws.Range("D60").HorizontalAlignment = -4131
$ws.Range("D60").NumberFormat = "0"
$ws.Range("E60").HorizontalAlignment = -4131
$ws.Range("E60").NumberFormat = "0"

# Update dimension / view: sheet scrolled back to top, selection moved to E3
$ws.Range("E3").Select()
